$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H, formatted like the other header cells (E1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "prop.ripe"

# New "prop.ripe" column: ripe / (mature + ripe + spent) for each data row
$ws.Range("H2").Formula = "=F2/SUM(E2:G2)"
$ws.Range("H3:H7").Formula = "=F3/SUM(E3:G3)"

# Match the author's final selection
[void]$ws.Range("H9").Select()
